$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Sheet1 cell values ---
$ws1.Range("G1").Value = -179.74
$ws1.Range("H1").Value = -179.74
$ws1.Range("I1").Value = -165.53
$ws1.Range("J1").Value = -165.53
$ws1.Range("K1").Value = -344.08
$ws1.Range("L1").Value = -344.08
$ws1.Range("M1").Value = -325.55
$ws1.Range("N1").Value = -325.55
$ws1.Range("O1").Value = -16.33
$ws1.Range("P1").Value = -16.33
$ws1.Range("Q1").Value = -3.77
$ws1.Range("R1").Value = -3.77
$ws1.Range("G2").Value = 1803.4
$ws1.Range("H2").Value = -1803.4
$ws1.Range("I2").Value = 1705.36
$ws1.Range("J2").Value = -1705.36
$ws1.Range("K2").Value = 240.83
$ws1.Range("L2").Value = -240.83
$ws1.Range("M2").Value = 226.63
$ws1.Range("N2").Value = -226.63
$ws1.Range("O2").Value = 36.02
$ws1.Range("P2").Value = -36.02
$ws1.Range("Q2").Value = 11.47
$ws1.Range("R2").Value = -11.47
$ws1.Range("G7").Value = 0.00966
$ws1.Range("H7").Value = 0.00966
$ws1.Range("I7").Value = 0.01138
$ws1.Range("J7").Value = 0.01138
$ws1.Range("K7").Value = 0.29783
$ws1.Range("L7").Value = 0.29783
$ws1.Range("M7").Value = 0.33335
$ws1.Range("N7").Value = 0.33335
$ws1.Range("O7").Value = 0.13866
$ws1.Range("P7").Value = 0.13866
$ws1.Range("Q7").Value = 0.02008
$ws1.Range("R7").Value = 0.02008
$ws1.Range("G8").Value = 0.00966
$ws1.Range("H8").Value = 0.00966
$ws1.Range("I8").Value = 0.01138
$ws1.Range("J8").Value = 0.01138
$ws1.Range("K8").Value = 0.2978
$ws1.Range("L8").Value = 0.2978
$ws1.Range("M8").Value = 0.33339
$ws1.Range("N8").Value = 0.33339
$ws1.Range("O8").Value = 0.13878
$ws1.Range("P8").Value = 0.13878
$ws1.Range("Q8").Value = 0.02012
$ws1.Range("R8").Value = 0.02012
$ws1.Range("G9").Value = 0.04903
$ws1.Range("H9").Value = 0.04903
$ws1.Range("I9").Value = 0.04765
$ws1.Range("J9").Value = 0.04765
$ws1.Range("K9").Value = 0.43976
$ws1.Range("L9").Value = 0.43976
$ws1.Range("M9").Value = 0.43029
$ws1.Range("N9").Value = 0.43029
$ws1.Range("O9").Value = 0.00159
$ws1.Range("P9").Value = 0.00159
$ws1.Range("Q9").Value = 0.00002
$ws1.Range("R9").Value = 0.00002
$ws1.Range("G10").Value = 0.04903
$ws1.Range("H10").Value = 0.04903
$ws1.Range("I10").Value = 0.04765
$ws1.Range("J10").Value = 0.04765
$ws1.Range("K10").Value = 0.43973
$ws1.Range("L10").Value = 0.43973
$ws1.Range("M10").Value = 0.43033
$ws1.Range("N10").Value = 0.43033
$ws1.Range("O10").Value = 0.00167
$ws1.Range("P10").Value = 0.00167
$ws1.Range("Q10").Value = 0.00013
$ws1.Range("R10").Value = 0.00013
$ws1.Range("G11").Value = 0.23589
$ws1.Range("H11").Value = 0.23589
$ws1.Range("I11").Value = 0.23833
$ws1.Range("J11").Value = 0.23833
$ws1.Range("K11").Value = 0.15563
$ws1.Range("L11").Value = 0.15563
$ws1.Range("M11").Value = 0.13471
$ws1.Range("N11").Value = 0.13471
$ws1.Range("O11").Value = 0.02137
$ws1.Range("P11").Value = 0.02137
$ws1.Range("Q11").Value = 0.00097
$ws1.Range("R11").Value = 0.00097
$ws1.Range("G12").Value = 0.23589
$ws1.Range("H12").Value = 0.23589
$ws1.Range("I12").Value = 0.23833
$ws1.Range("J12").Value = 0.23833
$ws1.Range("K12").Value = 0.15563
$ws1.Range("L12").Value = 0.15563
$ws1.Range("M12").Value = 0.13471
$ws1.Range("N12").Value = 0.13471
$ws1.Range("O12").Value = 0.02137
$ws1.Range("P12").Value = 0.02137
$ws1.Range("Q12").Value = 0.00097
$ws1.Range("R12").Value = 0.00097
$ws1.Range("G13").Value = 0.24375
$ws1.Range("H13").Value = 0.24375
$ws1.Range("I13").Value = 0.24202
$ws1.Range("J13").Value = 0.24202
$ws1.Range("K13").Value = 0.01471
$ws1.Range("L13").Value = 0.01471
$ws1.Range("M13").Value = 0.01256
$ws1.Range("N13").Value = 0.01256
$ws1.Range("O13").Value = 0.00083
$ws1.Range("P13").Value = 0.00083
$ws1.Range("Q13").Value = 0.00004
$ws1.Range("R13").Value = 0.00004
$ws1.Range("G14").Value = 0.24375
$ws1.Range("H14").Value = 0.24375
$ws1.Range("I14").Value = 0.24203
$ws1.Range("J14").Value = 0.24203
$ws1.Range("K14").Value = 0.01471
$ws1.Range("L14").Value = 0.01471
$ws1.Range("M14").Value = 0.01256
$ws1.Range("N14").Value = 0.01256
$ws1.Range("O14").Value = 0.00079
$ws1.Range("P14").Value = 0.00079
$ws1.Range("Q14").Value = 0.00004
$ws1.Range("R14").Value = 0.00004
$ws1.Range("G15").Value = 0.04244
$ws1.Range("H15").Value = 0.04244
$ws1.Range("I15").Value = 0.0423
$ws1.Range("J15").Value = 0.0423
$ws1.Range("K15").Value = 0.04926
$ws1.Range("L15").Value = 0.04926
$ws1.Range("M15").Value = 0.03088
$ws1.Range("N15").Value = 0.03088
$ws1.Range("O15").Value = 0.40147
$ws1.Range("P15").Value = 0.40147
$ws1.Range("Q15").Value = 0.06171
$ws1.Range("R15").Value = 0.06171
$ws1.Range("G16").Value = 0.04244
$ws1.Range("H16").Value = 0.04244
$ws1.Range("I16").Value = 0.0423
$ws1.Range("J16").Value = 0.0423
$ws1.Range("K16").Value = 0.04926
$ws1.Range("L16").Value = 0.04926
$ws1.Range("M16").Value = 0.03088
$ws1.Range("N16").Value = 0.03088
$ws1.Range("O16").Value = 0.40147
$ws1.Range("P16").Value = 0.40147
$ws1.Range("Q16").Value = 0.06171
$ws1.Range("R16").Value = 0.06171
$ws1.Range("O17").Value = 0.06374
$ws1.Range("P17").Value = 0.06374
$ws1.Range("Q17").Value = 0.56181
$ws1.Range("R17").Value = 0.56181
$ws1.Range("O18").Value = 0.0589
$ws1.Range("P18").Value = 0.0589
$ws1.Range("Q18").Value = 0.60738
$ws1.Range("R18").Value = 0.60738

# --- Add Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 header rows (mirrors updated Sheet1 rows 1-2) ---
$ws2.Range("A1").Value = "Real part"
$ws2.Range("C1").Value = -393.71
$ws2.Range("D1").Value = -154.09
$ws2.Range("E1").Value = -15.89
$ws2.Range("F1").Value = -15.89
$ws2.Range("G1").Value = -179.74
$ws2.Range("H1").Value = -179.74
$ws2.Range("I1").Value = -165.53
$ws2.Range("J1").Value = -165.53
$ws2.Range("K1").Value = -344.08
$ws2.Range("L1").Value = -344.08
$ws2.Range("M1").Value = -325.55
$ws2.Range("N1").Value = -325.55
$ws2.Range("O1").Value = -16.33
$ws2.Range("P1").Value = -16.33
$ws2.Range("Q1").Value = -3.77
$ws2.Range("R1").Value = -3.77
$ws2.Range("A2").Value = "Image part"
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 12.24
$ws2.Range("F2").Value = -12.24
$ws2.Range("G2").Value = 1803.4
$ws2.Range("H2").Value = -1803.4
$ws2.Range("I2").Value = 1705.36
$ws2.Range("J2").Value = -1705.36
$ws2.Range("K2").Value = 240.83
$ws2.Range("L2").Value = -240.83
$ws2.Range("M2").Value = 226.63
$ws2.Range("N2").Value = -226.63
$ws2.Range("O2").Value = 36.02
$ws2.Range("P2").Value = -36.02
$ws2.Range("Q2").Value = 11.47
$ws2.Range("R2").Value = -11.47

# --- Populate Sheet2 new parameter rows ---
$ws2.Range("B3").Value = "Rg"
$ws2.Range("B4").Value = "Lg"

# --- Activate Sheet2 and set its selection ---
$ws2.Activate()
$ws2.Range("C2:R2").Select()
